$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# Locate the "Date" column by scanning the header row (row 1) instead of
# hard-coding a column letter.
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "Date") {
        $dateCol = $c
        break
    }
}

if ($dateCol -gt 0) {
    # The bad values were shifted one day early (e.g. "6-15-2012-13" for
    # what should be the ISO date 2013-06-15). Force the column to text
    # formatting first so Excel doesn't reinterpret the corrected,
    # ISO-looking date string as a date serial number.
    $colRange = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
    $colRange.NumberFormat = "@"

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $dateCol)
        if ($cell.Text -eq "6-15-2012-13") {
            $cell.Value = "2013-06-15"
        }
    }
}
